# Update betting-odds values on the active worksheet (flashscore weekly odds sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.4
$ws.Range("H2").Value = 4.2
$ws.Range("I2").Value = 9
$ws.Range("L2").Value = 8.5
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("AA2").Value = 2.5
$ws.Range("AB2").Value = 1.5
$ws.Range("AC2").Value = 5
$ws.Range("AF2").Value = 8.5
$ws.Range("AM2").Value = 17
$ws.Range("AR2").Value = 81

# Row 3
$ws.Range("N3").Value = 8
$ws.Range("U3").Value = 3.7
$ws.Range("V3").Value = 1.27

# Row 4
$ws.Range("G4").Value = 2.63
$ws.Range("I4").Value = 2.63
$ws.Range("L4").Value = 3.5
$ws.Range("S4").Value = 2.25
$ws.Range("T4").Value = 1.62
$ws.Range("W4").Value = 4
$ws.Range("X4").Value = 1.22
$ws.Range("AC4").Value = 7.5
$ws.Range("AD4").Value = 12
$ws.Range("AF4").Value = 26
$ws.Range("AK4").Value = 17
$ws.Range("AO4").Value = 11

# Row 5
$ws.Range("G5").Value = 1.95
$ws.Range("I5").Value = 3.6
$ws.Range("J5").Value = 2.63
$ws.Range("L5").Value = 4.33
$ws.Range("W5").Value = 3.5
$ws.Range("X5").Value = 1.29
$ws.Range("AD5").Value = 9
$ws.Range("AF5").Value = 17
$ws.Range("AK5").Value = 17
$ws.Range("AM5").Value = 10
$ws.Range("AN5").Value = 19
$ws.Range("AO5").Value = 13
